$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, reflecting the reordering of the
# weekly Fruta/Hortaliza records (dates + their associated Volumen,
# Precio minimo/maximo/promedio, Precio $/Kg columns).
$updates = @{
    2  = @{ D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 }
    3  = @{ D = 44902; M = 35; N = 12000; O = 12000; P = 12000; S = 2400 }
    4  = @{ D = 44914; M = 56; N = 23000; O = 23000; P = 23000; S = 4600 }
    5  = @{ D = 44188; M = 30 }
    6  = @{ D = 44193; M = 40 }
    7  = @{ D = 44931; M = 50; N = 18000; O = 18000; P = 18000; S = 3600 }
    9  = @{ D = 44196; M = 56 }
    10 = @{ D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 }
    11 = @{ D = 44186; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    12 = @{ D = 44907; M = 45; N = 25000; O = 25000; P = 25000; S = 5000 }
    13 = @{ D = 44189; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
